$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 457.14285
$ws.Range("I8").Value = 457.14285
$ws.Range("K8").Value = 1371.42855
$ws.Range("M8").Value = -1232.42855
$ws.Range("H43").Value = 7951579.5
$ws.Range("J43").Value = 11111611
$ws.Range("L43").Value = 11111611
$ws.Range("N43").Value = -11111749
$ws.Range("H107").Value = 1778.45
$ws.Range("I107").Value = 1386.0588
$ws.Range("K107").Value = 1386.0588
$ws.Range("M107").Value = 533.9412
$ws.Range("H132").Value = 12353409
$ws.Range("I132").Value = 15879197
$ws.Range("J132").Value = 13152
$ws.Range("K132").Value = 47637591
$ws.Range("L132").Value = 39456
$ws.Range("M132").Value = -47635061
$ws.Range("N132").Value = -44516
$ws.Range("H136").Value = 26473.334
$ws.Range("J136").Value = 28768
$ws.Range("L136").Value = 28768
$ws.Range("N136").Value = -38968

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1324.35
$ws.Range("I132").Value = 1002.2353
$ws.Range("K132").Value = 3006.7059
$ws.Range("M132").Value = -476.7058999999999

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 582.625
$ws.Range("I22").Value = 443.66666
$ws.Range("J22").Value = 999.5
$ws.Range("K22").Value = 443.66666
$ws.Range("L22").Value = 999.5
$ws.Range("M22").Value = -270.66666
$ws.Range("N22").Value = -1345.5
$ws.Range("H134").Value = 1406.6123
$ws.Range("I134").Value = 929.0513
$ws.Range("K134").Value = 2787.1539
$ws.Range("M134").Value = -252.1538999999998

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 16666.666
$ws.Range("I56").Value = 16000
$ws.Range("J56").Value = 17000
$ws.Range("K56").Value = 16000
$ws.Range("L56").Value = 17000
$ws.Range("M56").Value = -15155
$ws.Range("N56").Value = -18690
$ws.Range("H86").Value = 3937158.2
$ws.Range("I86").Value = 8336946
$ws.Range("J86").Value = 26235.777
$ws.Range("K86").Value = 8336946
$ws.Range("L86").Value = 26235.777
$ws.Range("M86").Value = -8335823
$ws.Range("N86").Value = -28481.777
$ws.Range("H89").Value = 3937158.2
$ws.Range("I89").Value = 8336946
$ws.Range("J89").Value = 26235.777
$ws.Range("K89").Value = 41684730
$ws.Range("L89").Value = 131178.885
$ws.Range("M89").Value = -41679114
$ws.Range("N89").Value = -142410.885
$ws.Range("H132").Value = 2849.8
$ws.Range("I132").Value = 1233.3334
$ws.Range("K132").Value = 3700.0002
$ws.Range("M132").Value = -1170.0002

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 69.48148
$ws.Range("I12").Value = 85.28570999999999
$ws.Range("J12").Value = 63.95
$ws.Range("K12").Value = 255.85713
$ws.Range("L12").Value = 191.85
$ws.Range("M12").Value = -82.85712999999998
$ws.Range("N12").Value = -537.85
$ws.Range("H64").Value = 3778.5518
$ws.Range("I64").Value = 1990
$ws.Range("J64").Value = 3911.037
$ws.Range("K64").Value = 5970
$ws.Range("L64").Value = 11733.111
$ws.Range("M64").Value = -5700
$ws.Range("N64").Value = -12273.111
$ws.Range("H67").Value = 3778.5518
$ws.Range("I67").Value = 1990
$ws.Range("J67").Value = 3911.037
$ws.Range("K67").Value = 5970
$ws.Range("L67").Value = 11733.111
$ws.Range("M67").Value = -5034
$ws.Range("N67").Value = -13605.111
$ws.Range("H113").Value = 692.1429000000001
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 692.1429000000001
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").Value = 2076.4287
$ws.Range("N113").Value = -6416.4287

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 11642.75
$ws.Range("I5").Value = 4416.6665
$ws.Range("J5").Value = 12917.941
$ws.Range("K5").Value = 4416.6665
$ws.Range("L5").Value = 12917.941
$ws.Range("M5").Value = -4304.6665
$ws.Range("N5").Value = -13141.941
$ws.Range("H80").Value = 6900
$ws.Range("I80").Value = 7000
$ws.Range("J80").Value = 6800
$ws.Range("K80").Value = 7000
$ws.Range("L80").Value = 6800
$ws.Range("M80").Value = -6002
$ws.Range("N80").Value = -8796
$ws.Range("H83").Value = 6900
$ws.Range("I83").Value = 7000
$ws.Range("J83").Value = 6800
$ws.Range("K83").Value = 35000
$ws.Range("L83").Value = 34000
$ws.Range("M83").Value = -30008
$ws.Range("N83").Value = -43984
$ws.Range("H88").Value = 40000
$ws.Range("J88").Value = 40000
$ws.Range("L88").Value = 40000
$ws.Range("N88").Value = -40902
$ws.Range("H91").Value = 40000
$ws.Range("J91").Value = 40000
$ws.Range("L91").Value = 40000
$ws.Range("N91").Value = -43120
$ws.Range("H104").Value = 48500
$ws.Range("J104").Value = 48500
$ws.Range("L104").Value = 48500
$ws.Range("N104").Value = -55488
$ws.Range("H107").Value = 941.5625
$ws.Range("I107").Value = 960.5454999999999
$ws.Range("K107").Value = 960.5454999999999
$ws.Range("M107").Value = 959.4545000000001

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1983.5
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1983.5
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").Value = 1983.5
$ws.Range("N22").Value = -2573.5
$ws.Range("H27").Value = 1983.5
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1983.5
$ws.Range("K27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("M27").Value = 1983.5
$ws.Range("N27").Value = -2197.5
$ws.Range("H87").Value = 22000
$ws.Range("J87").Value = 22000
$ws.Range("L87").Value = 22000
$ws.Range("N87").Value = -24246
$ws.Range("H90").Value = 22000
$ws.Range("J90").Value = 22000
$ws.Range("L90").Value = 66000
$ws.Range("N90").Value = -77232
$ws.Range("H132").Value = 27293.564
$ws.Range("I132").Value = 1049.826
$ws.Range("J132").Value = 65018.938
$ws.Range("K132").Value = 3149.478
$ws.Range("L132").Value = 195056.814
$ws.Range("M132").Value = -619.4780000000001
$ws.Range("N132").Value = -200116.814

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 17900
$ws.Range("J64").Value = 17900
$ws.Range("L64").Value = 17900
$ws.Range("N64").Value = -18396
$ws.Range("H67").Value = 17900
$ws.Range("J67").Value = 17900
$ws.Range("L67").Value = 17900
$ws.Range("N67").Value = -19616
$ws.Range("H107").Value = 586.9
$ws.Range("I107").Value = 481.2857
$ws.Range("J107").Value = 833.3333
$ws.Range("K107").Value = 1443.8571
$ws.Range("L107").Value = 2499.9999
$ws.Range("M107").Value = 476.1428999999998
$ws.Range("N107").Value = -6339.9999
